# Apply the data refresh: rows 2-8 (the species observation records) were
# re-pulled from source and ended up in a different row order, with each
# row's Id/TaxonId/name/author/time fields changing to match.  Only columns
# A, B, D, E, F, G, H, Z, AB differ between the old and new row order; every
# other column (location, date, observer, etc.) is identical across all
# rows, so we only need to rewrite those nine columns per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
    2 = @{ A = 111708920; B = 90666; D = "LC"; E = 4364; F = "Dropptaggsvamp";      G = "Hydnellum ferrugineum"; H = "(Fr.:Fr.) P. Karst.";      Z = "13:53" }
    3 = @{ A = 111708099; B = 90660; D = "NT"; E = 4362; F = "Blå taggsvamp";       G = "Hydnellum caeruleum";   H = "(Hornem.) P.Karst.";        Z = "14:16" }
    4 = @{ A = 111708029; B = 90662; D = "LC"; E = 4363; F = "Zontaggsvamp";        G = "Hydnellum concrescens"; H = "(Pers.) Banker";             Z = "14:21" }
    5 = @{ A = 111708162; B = 90658; D = "NT"; E = 4361; F = "Orange taggsvamp";    G = "Hydnellum aurantiacum"; H = "(Batsch:Fr.) P.Karst.";      Z = "14:12" }
    6 = @{ A = 111706580; B = 88032; D = "VU"; E = 6276; F = "Goliatmusseron";      G = "Tricholoma matsutake";  H = "(S.Ito & S.Imai) Singer";    Z = "14:48" }
    7 = @{ A = 111708888; B = 90678; D = "LC"; E = 4366; F = "Skarp dropptaggsvamp"; G = "Hydnellum peckii";     H = "Banker";                     Z = "13:54" }
    8 = @{ A = 111704319; B = 90710; D = "NT"; E = 5449; F = "Svart taggsvamp";     G = "Phellodon niger";       H = "(Fr.:Fr.) P.Karst.";          Z = "15:11" }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 1).Value  = $vals.A   # A: Id
    $ws.Cells.Item($r, 2).Value  = $vals.B   # B: Taxonsorteringsordning
    $ws.Cells.Item($r, 4).Value  = $vals.D   # D: Rödlistade
    $ws.Cells.Item($r, 5).Value  = $vals.E   # E: TaxonId
    $ws.Cells.Item($r, 6).Value  = $vals.F   # F: Artnamn
    $ws.Cells.Item($r, 7).Value  = $vals.G   # G: Vetenskapligt namn
    $ws.Cells.Item($r, 8).Value  = $vals.H   # H: Auktor
    $ws.Cells.Item($r, 26).Value = $vals.Z   # Z: Starttid
    $ws.Cells.Item($r, 28).Value = $vals.Z   # AB: Sluttid
}
